$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").Value = 2019
$ws.Range("B62").Value = 0
$ws.Range("B62").NumberFormat = $ws.Range("B61").NumberFormat

$ws.Range("A63").Value = 2020
$ws.Range("B63").Value = 0
$ws.Range("B63").NumberFormat = $ws.Range("B61").NumberFormat

$ws.Range("K56").Select()
